# Add a new worksheet "Dire Rat AC" as a copy of the existing "Erdlu_AC"
# sheet (same layout/formulas), then update the Feat (E) and Feat-number
# (H) columns, and the starting Companion HD (G6), to describe the Dire
# Rat animal companion instead of the Erdlu one.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the Erdlu_AC sheet right after itself, then rename the copy.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Dire Rat AC"

# --- First progression table (levels 1st-2nd .. 18th-20th) ---

# New strings are introduced in this order so the saved shared-string
# table lines up with the source workbook (Alterness, WF(Creature),
# Spr Attack, Iron Will, EWF (Creature), Epic Prowess).
$ws2.Range("E6").Value = "Alterness"
$ws2.Range("G6").Value = 1

$ws2.Range("E7").Value = "Dodge"

$ws2.Range("E8").Value = "-"
$ws2.Range("H8").Value = "-"

$ws2.Range("E9").Value = "WF(Creature)"
$ws2.Range("H9").Value = "Feat 3"

$ws2.Range("E10").Value = "Mobility"

$ws2.Range("E11").Value = "-"
$ws2.Range("H11").Value = "-"

$ws2.Range("E12").Value = "Blind-fight"
$ws2.Range("H12").Value = "Feat 5"

# --- Epic progression table (levels 21st-23rd .. 39th-40th) ---

$ws2.Range("E16").Value = "Spr Attack"

$ws2.Range("E17").Value = "-"
$ws2.Range("H17").Value = "-"

$ws2.Range("E18").Value = "Iron Will"
$ws2.Range("H18").Value = "Feat 7"

$ws2.Range("E20").Value = "-"
$ws2.Range("H20").Value = "-"

$ws2.Range("E21").Value = "EWF (Creature)"

$ws2.Range("E19").Value = "Epic Prowess"

# Match the saved selection on the new sheet.
$ws2.Range("E22").Select()
